# ---------------------------------------------------------------
# 1) Insert a new worksheet '2022-Q3' right before '2022-Q2', shifting
#    2022-Q2 .. 2021-Q1 one slot to the right (same order, same content).
# ---------------------------------------------------------------
$wb = $excel.ActiveWorkbook
$insertBefore = $wb.Worksheets.Item("2022-Q2")
$newSheet = $wb.Worksheets.Add($insertBefore)
$newSheet.Name = "2022-Q3"

# Re-fetch by name: the Add(Before) reference aliases the same slot as $newSheet
# until re-looked-up, so grab a fresh handle on "2022-Q2" to use as a format donor.
$q2 = $wb.Worksheets.Item("2022-Q2")

# ---------------------------------------------------------------
# 2) Populate the new '2022-Q3' sheet with its header + 16 data rows.
#    Values that look numeric ('005805', '4.09', '93.86', ...) are written
#    with a leading apostrophe + ClearFormats() so they land as plain text
#    cells (matching source columns B-G), not auto-coerced numbers.
# ---------------------------------------------------------------
$newSheet.Cells.Item(1,2).Value = '''基金代码'
$newSheet.Cells.Item(1,2).ClearFormats()
$newSheet.Cells.Item(1,3).Value = '''基金名称'
$newSheet.Cells.Item(1,3).ClearFormats()
$newSheet.Cells.Item(1,4).Value = '''基金规模'
$newSheet.Cells.Item(1,4).ClearFormats()
$newSheet.Cells.Item(1,5).Value = '''股票总仓位'
$newSheet.Cells.Item(1,5).ClearFormats()
$newSheet.Cells.Item(1,6).Value = '''仓位占比'
$newSheet.Cells.Item(1,6).ClearFormats()
$newSheet.Cells.Item(1,7).Value = '''持有市值(亿元)'
$newSheet.Cells.Item(1,7).ClearFormats()
$newSheet.Cells.Item(1,8).Value = '''仓位排名'
$newSheet.Cells.Item(1,8).ClearFormats()

$newSheet.Cells.Item(2,1).Value = 0
$newSheet.Cells.Item(2,1).ClearFormats()
$newSheet.Cells.Item(2,2).Value = '''005805'
$newSheet.Cells.Item(2,2).ClearFormats()
$newSheet.Cells.Item(2,3).Value = '''华泰柏瑞医疗健康混合A'
$newSheet.Cells.Item(2,3).ClearFormats()
$newSheet.Cells.Item(2,4).Value = '''4.09'
$newSheet.Cells.Item(2,4).ClearFormats()
$newSheet.Cells.Item(2,5).Value = '''93.86'
$newSheet.Cells.Item(2,5).ClearFormats()
$newSheet.Cells.Item(2,6).Value = '''4.45'
$newSheet.Cells.Item(2,6).ClearFormats()
$newSheet.Cells.Item(2,7).Value = '''0.1820'
$newSheet.Cells.Item(2,7).ClearFormats()
$newSheet.Cells.Item(2,8).Value = 6
$newSheet.Cells.Item(2,8).ClearFormats()

$newSheet.Cells.Item(3,1).Value = 1
$newSheet.Cells.Item(3,1).ClearFormats()
$newSheet.Cells.Item(3,2).Value = '''009877'
$newSheet.Cells.Item(3,2).ClearFormats()
$newSheet.Cells.Item(3,3).Value = '''中银内核驱动股票A'
$newSheet.Cells.Item(3,3).ClearFormats()
$newSheet.Cells.Item(3,4).Value = '''2.44'
$newSheet.Cells.Item(3,4).ClearFormats()
$newSheet.Cells.Item(3,5).Value = '''90.00'
$newSheet.Cells.Item(3,5).ClearFormats()
$newSheet.Cells.Item(3,6).Value = '''4.79'
$newSheet.Cells.Item(3,6).ClearFormats()
$newSheet.Cells.Item(3,7).Value = '''0.1169'
$newSheet.Cells.Item(3,7).ClearFormats()
$newSheet.Cells.Item(3,8).Value = 10
$newSheet.Cells.Item(3,8).ClearFormats()

$newSheet.Cells.Item(4,1).Value = 2
$newSheet.Cells.Item(4,1).ClearFormats()
$newSheet.Cells.Item(4,2).Value = '''470888'
$newSheet.Cells.Item(4,2).ClearFormats()
$newSheet.Cells.Item(4,3).Value = '''汇添富香港优势精选混合（QDII）'
$newSheet.Cells.Item(4,3).ClearFormats()
$newSheet.Cells.Item(4,4).Value = '''1.63'
$newSheet.Cells.Item(4,4).ClearFormats()
$newSheet.Cells.Item(4,5).Value = '''78.50'
$newSheet.Cells.Item(4,5).ClearFormats()
$newSheet.Cells.Item(4,6).Value = '''7.06'
$newSheet.Cells.Item(4,6).ClearFormats()
$newSheet.Cells.Item(4,7).Value = '''0.1151'
$newSheet.Cells.Item(4,7).ClearFormats()
$newSheet.Cells.Item(4,8).Value = 2
$newSheet.Cells.Item(4,8).ClearFormats()

$newSheet.Cells.Item(5,1).Value = 3
$newSheet.Cells.Item(5,1).ClearFormats()
$newSheet.Cells.Item(5,2).Value = '''007718'
$newSheet.Cells.Item(5,2).ClearFormats()
$newSheet.Cells.Item(5,3).Value = '''中银创新医疗混合A'
$newSheet.Cells.Item(5,3).ClearFormats()
$newSheet.Cells.Item(5,4).Value = '''2.82'
$newSheet.Cells.Item(5,4).ClearFormats()
$newSheet.Cells.Item(5,5).Value = '''92.35'
$newSheet.Cells.Item(5,5).ClearFormats()
$newSheet.Cells.Item(5,6).Value = '''3.77'
$newSheet.Cells.Item(5,6).ClearFormats()
$newSheet.Cells.Item(5,7).Value = '''0.1063'
$newSheet.Cells.Item(5,7).ClearFormats()
$newSheet.Cells.Item(5,8).Value = 7
$newSheet.Cells.Item(5,8).ClearFormats()

$newSheet.Cells.Item(6,1).Value = 4
$newSheet.Cells.Item(6,1).ClearFormats()
$newSheet.Cells.Item(6,2).Value = '''012584'
$newSheet.Cells.Item(6,2).ClearFormats()
$newSheet.Cells.Item(6,3).Value = '''南方中国新兴经济9个月持有期混合（QDII）A'
$newSheet.Cells.Item(6,3).ClearFormats()
$newSheet.Cells.Item(6,4).Value = '''2.69'
$newSheet.Cells.Item(6,4).ClearFormats()
$newSheet.Cells.Item(6,5).Value = '''91.51'
$newSheet.Cells.Item(6,5).ClearFormats()
$newSheet.Cells.Item(6,6).Value = '''3.87'
$newSheet.Cells.Item(6,6).ClearFormats()
$newSheet.Cells.Item(6,7).Value = '''0.1041'
$newSheet.Cells.Item(6,7).ClearFormats()
$newSheet.Cells.Item(6,8).Value = 5
$newSheet.Cells.Item(6,8).ClearFormats()

$newSheet.Cells.Item(7,1).Value = 5
$newSheet.Cells.Item(7,1).ClearFormats()
$newSheet.Cells.Item(7,2).Value = '''513120'
$newSheet.Cells.Item(7,2).ClearFormats()
$newSheet.Cells.Item(7,3).Value = '''广发中证香港创新药（QDII-ETF）'
$newSheet.Cells.Item(7,3).ClearFormats()
$newSheet.Cells.Item(7,4).Value = '''1.09'
$newSheet.Cells.Item(7,4).ClearFormats()
$newSheet.Cells.Item(7,5).Value = '''98.58'
$newSheet.Cells.Item(7,5).ClearFormats()
$newSheet.Cells.Item(7,6).Value = '''2.84'
$newSheet.Cells.Item(7,6).ClearFormats()
$newSheet.Cells.Item(7,7).Value = '''0.0310'
$newSheet.Cells.Item(7,7).ClearFormats()
$newSheet.Cells.Item(7,8).Value = 10
$newSheet.Cells.Item(7,8).ClearFormats()

$newSheet.Cells.Item(8,1).Value = 6
$newSheet.Cells.Item(8,1).ClearFormats()
$newSheet.Cells.Item(8,2).Value = '''005029'
$newSheet.Cells.Item(8,2).ClearFormats()
$newSheet.Cells.Item(8,3).Value = '''中银产业精选混合'
$newSheet.Cells.Item(8,3).ClearFormats()
$newSheet.Cells.Item(8,4).Value = '''0.44'
$newSheet.Cells.Item(8,4).ClearFormats()
$newSheet.Cells.Item(8,5).Value = '''83.28'
$newSheet.Cells.Item(8,5).ClearFormats()
$newSheet.Cells.Item(8,6).Value = '''6.04'
$newSheet.Cells.Item(8,6).ClearFormats()
$newSheet.Cells.Item(8,7).Value = '''0.0266'
$newSheet.Cells.Item(8,7).ClearFormats()
$newSheet.Cells.Item(8,8).Value = 7
$newSheet.Cells.Item(8,8).ClearFormats()

$newSheet.Cells.Item(9,1).Value = 7
$newSheet.Cells.Item(9,1).ClearFormats()
$newSheet.Cells.Item(9,2).Value = '''011453'
$newSheet.Cells.Item(9,2).ClearFormats()
$newSheet.Cells.Item(9,3).Value = '''华泰柏瑞医疗健康混合C'
$newSheet.Cells.Item(9,3).ClearFormats()
$newSheet.Cells.Item(9,4).Value = '''0.43'
$newSheet.Cells.Item(9,4).ClearFormats()
$newSheet.Cells.Item(9,5).Value = '''93.86'
$newSheet.Cells.Item(9,5).ClearFormats()
$newSheet.Cells.Item(9,6).Value = '''4.45'
$newSheet.Cells.Item(9,6).ClearFormats()
$newSheet.Cells.Item(9,7).Value = '''0.0191'
$newSheet.Cells.Item(9,7).ClearFormats()
$newSheet.Cells.Item(9,8).Value = 6
$newSheet.Cells.Item(9,8).ClearFormats()

$newSheet.Cells.Item(10,1).Value = 8
$newSheet.Cells.Item(10,1).ClearFormats()
$newSheet.Cells.Item(10,2).Value = '''010783'
$newSheet.Cells.Item(10,2).ClearFormats()
$newSheet.Cells.Item(10,3).Value = '''德邦沪港深龙头混合A'
$newSheet.Cells.Item(10,3).ClearFormats()
$newSheet.Cells.Item(10,4).Value = '''0.55'
$newSheet.Cells.Item(10,4).ClearFormats()
$newSheet.Cells.Item(10,5).Value = '''84.96'
$newSheet.Cells.Item(10,5).ClearFormats()
$newSheet.Cells.Item(10,6).Value = '''3.27'
$newSheet.Cells.Item(10,6).ClearFormats()
$newSheet.Cells.Item(10,7).Value = '''0.0180'
$newSheet.Cells.Item(10,7).ClearFormats()
$newSheet.Cells.Item(10,8).Value = 7
$newSheet.Cells.Item(10,8).ClearFormats()

$newSheet.Cells.Item(11,1).Value = 9
$newSheet.Cells.Item(11,1).ClearFormats()
$newSheet.Cells.Item(11,2).Value = '''013897'
$newSheet.Cells.Item(11,2).ClearFormats()
$newSheet.Cells.Item(11,3).Value = '''德邦港股通成长精选混合型证券投资基金A'
$newSheet.Cells.Item(11,3).ClearFormats()
$newSheet.Cells.Item(11,4).Value = '''0.41'
$newSheet.Cells.Item(11,4).ClearFormats()
$newSheet.Cells.Item(11,5).Value = '''79.99'
$newSheet.Cells.Item(11,5).ClearFormats()
$newSheet.Cells.Item(11,6).Value = '''3.34'
$newSheet.Cells.Item(11,6).ClearFormats()
$newSheet.Cells.Item(11,7).Value = '''0.0137'
$newSheet.Cells.Item(11,7).ClearFormats()
$newSheet.Cells.Item(11,8).Value = 8
$newSheet.Cells.Item(11,8).ClearFormats()

$newSheet.Cells.Item(12,1).Value = 10
$newSheet.Cells.Item(12,1).ClearFormats()
$newSheet.Cells.Item(12,2).Value = '''013898'
$newSheet.Cells.Item(12,2).ClearFormats()
$newSheet.Cells.Item(12,3).Value = '''德邦港股通成长精选混合型证券投资基金C'
$newSheet.Cells.Item(12,3).ClearFormats()
$newSheet.Cells.Item(12,4).Value = '''0.37'
$newSheet.Cells.Item(12,4).ClearFormats()
$newSheet.Cells.Item(12,5).Value = '''79.99'
$newSheet.Cells.Item(12,5).ClearFormats()
$newSheet.Cells.Item(12,6).Value = '''3.34'
$newSheet.Cells.Item(12,6).ClearFormats()
$newSheet.Cells.Item(12,7).Value = '''0.0124'
$newSheet.Cells.Item(12,7).ClearFormats()
$newSheet.Cells.Item(12,8).Value = 8
$newSheet.Cells.Item(12,8).ClearFormats()

$newSheet.Cells.Item(13,1).Value = 11
$newSheet.Cells.Item(13,1).ClearFormats()
$newSheet.Cells.Item(13,2).Value = '''010784'
$newSheet.Cells.Item(13,2).ClearFormats()
$newSheet.Cells.Item(13,3).Value = '''德邦沪港深龙头混合C'
$newSheet.Cells.Item(13,3).ClearFormats()
$newSheet.Cells.Item(13,4).Value = '''0.36'
$newSheet.Cells.Item(13,4).ClearFormats()
$newSheet.Cells.Item(13,5).Value = '''84.96'
$newSheet.Cells.Item(13,5).ClearFormats()
$newSheet.Cells.Item(13,6).Value = '''3.27'
$newSheet.Cells.Item(13,6).ClearFormats()
$newSheet.Cells.Item(13,7).Value = '''0.0118'
$newSheet.Cells.Item(13,7).ClearFormats()
$newSheet.Cells.Item(13,8).Value = 7
$newSheet.Cells.Item(13,8).ClearFormats()

$newSheet.Cells.Item(14,1).Value = 12
$newSheet.Cells.Item(14,1).ClearFormats()
$newSheet.Cells.Item(14,2).Value = '''010500'
$newSheet.Cells.Item(14,2).ClearFormats()
$newSheet.Cells.Item(14,3).Value = '''中银创新医疗混合C'
$newSheet.Cells.Item(14,3).ClearFormats()
$newSheet.Cells.Item(14,4).Value = '''0.18'
$newSheet.Cells.Item(14,4).ClearFormats()
$newSheet.Cells.Item(14,5).Value = '''92.35'
$newSheet.Cells.Item(14,5).ClearFormats()
$newSheet.Cells.Item(14,6).Value = '''3.77'
$newSheet.Cells.Item(14,6).ClearFormats()
$newSheet.Cells.Item(14,7).Value = '''0.0068'
$newSheet.Cells.Item(14,7).ClearFormats()
$newSheet.Cells.Item(14,8).Value = 7
$newSheet.Cells.Item(14,8).ClearFormats()

$newSheet.Cells.Item(15,1).Value = 13
$newSheet.Cells.Item(15,1).ClearFormats()
$newSheet.Cells.Item(15,2).Value = '''006603'
$newSheet.Cells.Item(15,2).ClearFormats()
$newSheet.Cells.Item(15,3).Value = '''嘉实互融精选股票'
$newSheet.Cells.Item(15,3).ClearFormats()
$newSheet.Cells.Item(15,4).Value = '''0.12'
$newSheet.Cells.Item(15,4).ClearFormats()
$newSheet.Cells.Item(15,5).Value = '''82.85'
$newSheet.Cells.Item(15,5).ClearFormats()
$newSheet.Cells.Item(15,6).Value = '''4.94'
$newSheet.Cells.Item(15,6).ClearFormats()
$newSheet.Cells.Item(15,7).Value = '''0.0059'
$newSheet.Cells.Item(15,7).ClearFormats()
$newSheet.Cells.Item(15,8).Value = 3
$newSheet.Cells.Item(15,8).ClearFormats()

$newSheet.Cells.Item(16,1).Value = 14
$newSheet.Cells.Item(16,1).ClearFormats()
$newSheet.Cells.Item(16,2).Value = '''012585'
$newSheet.Cells.Item(16,2).ClearFormats()
$newSheet.Cells.Item(16,3).Value = '''南方中国新兴经济9个月持有期混合（QDII）C'
$newSheet.Cells.Item(16,3).ClearFormats()
$newSheet.Cells.Item(16,4).Value = '''0.10'
$newSheet.Cells.Item(16,4).ClearFormats()
$newSheet.Cells.Item(16,5).Value = '''91.51'
$newSheet.Cells.Item(16,5).ClearFormats()
$newSheet.Cells.Item(16,6).Value = '''3.87'
$newSheet.Cells.Item(16,6).ClearFormats()
$newSheet.Cells.Item(16,7).Value = '''0.0039'
$newSheet.Cells.Item(16,7).ClearFormats()
$newSheet.Cells.Item(16,8).Value = 5
$newSheet.Cells.Item(16,8).ClearFormats()

$newSheet.Cells.Item(17,1).Value = 15
$newSheet.Cells.Item(17,1).ClearFormats()
$newSheet.Cells.Item(17,2).Value = '''012600'
$newSheet.Cells.Item(17,2).ClearFormats()
$newSheet.Cells.Item(17,3).Value = '''中银内核驱动股票C'
$newSheet.Cells.Item(17,3).ClearFormats()
$newSheet.Cells.Item(17,4).Value = '''0.02'
$newSheet.Cells.Item(17,4).ClearFormats()
$newSheet.Cells.Item(17,5).Value = '''90.00'
$newSheet.Cells.Item(17,5).ClearFormats()
$newSheet.Cells.Item(17,6).Value = '''4.79'
$newSheet.Cells.Item(17,6).ClearFormats()
$newSheet.Cells.Item(17,7).Value = '''0.0010'
$newSheet.Cells.Item(17,7).ClearFormats()
$newSheet.Cells.Item(17,8).Value = 10
$newSheet.Cells.Item(17,8).ClearFormats()

# ---------------------------------------------------------------
# 3) Re-apply the header/index styling (cellXf s="2": bold, bordered,
#    centered) copied from the equivalent cells on the donor '2022-Q2' sheet.
#    Done LAST so ClearFormats() above never wipes it back out.
# ---------------------------------------------------------------
$q2.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)
$q2.Range("A2").Copy()
$newSheet.Range("A2:A17").PasteSpecial(-4122)
$newSheet.Range("A1").Select()

# ---------------------------------------------------------------
# 4) Update the summary ('总计') sheet: insert a new row 2 for 2022-Q3 and
#    bump the existing rows' running index (column A) down by one.
# ---------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")
$summary.Rows.Item(2).Insert()

$summary.Cells.Item(2,1).Value = 0
$summary.Cells.Item(2,2).Value = '''2022-Q3'
$summary.Cells.Item(2,3).Value = 16
$summary.Cells.Item(2,4).Value = 0.77
# Row-insert inherits the header row's bordered/bold style onto the new
# row's B:D cells; strip that back to plain (unstyled), matching the target.
$summary.Range("B2:D2").ClearFormats()

# Re-apply the column-A index style (s="2") to the freshly inserted row,
# copied from the row below (which already carries it).
$summary.Cells.Item(3,1).Copy()
$summary.Cells.Item(2,1).PasteSpecial(-4122)
$summary.Cells.Item(2,1).Value = 0

# Bump the running index in column A for every pre-existing row (now rows 3..7).
For ($r = 3; $r -le 7; $r++) {
    $summary.Cells.Item($r,1).Value = $r - 2
}

$summary.Range("A1").Select()
